$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Alobid%I%coreGivesNoEmail%1,            Arancibia%C%coreGivesNoEmail%1,            Avil\u00e9s-Jurado%FX%coreGivesNoEmail%1,            Bernal-Sprekelsen%M%coreGivesNoEmail%1,            Castro%P%coreGivesNoEmail%1,            de%Oss\u00f3 J%coreGivesNoEmail%1,            Gonz\u00e1lez-S\u00e1nchez%N%coreGivesNoEmail%1,            Guilemany%JM%coreGivesNoEmail%1,            Jes\u00fas%Rojas-Lechuga M%coreGivesNoEmail%1,            Langdon%C%coreGivesNoEmail%1,            Larrosa%F%coreGivesNoEmail%1,            Lehrer-Coriat%E%coreGivesNoEmail%1,            L\u00f3pez-Chac\u00f3n%M%coreGivesNoEmail%1,            Prieto%Alhambra D%coreGivesNoEmail%1,            Remacha%J%coreGivesNoEmail%1,            Ruiz-Sevilla%L%coreGivesNoEmail%1,            S\u00e1nchez%I%coreGivesNoEmail%1,            Vilaseca%I%coreGivesNoEmail%1]'
$ws.Range("E3").Value = '[Paul E.%Kwak%xref no email%1,        Joseph R.%Connors%xref no email%1,        Peter A.%Benedict%xref no email%1,        Micah R.%Timen%xref no email%1,        Binhuan%Wang%xref no email%1,        Yan%Zhang%xref no email%2,        Stephanie%Youlios%xref no email%1,        Kimberly%Sureau%xref no email%2,        Michael J.%Persky%xref no email%1,        Samaan%Rafeq%xref no email%2,        Luis%Angel%xref no email%1,        Milan R.%Amin%xref no email%2]'
$ws.Range("E4").Value = '[Kamil%Polok%NULL%1,                    Jakub%Fronczek%NULL%1,                    Peter Vernon%van Heerden%NULL%1,                    Hans%Flaatten%NULL%1,                    Bertrand%Guidet%NULL%1,                    Dylan W.%De Lange%NULL%1,                    Jesper%Fjølner%NULL%1,                    Susannah%Leaver%NULL%1,                    Michael%Beil%NULL%1,                    Sigal%Sviri%NULL%1,                    Raphael Romano%Bruno%NULL%1,                    Bernhard%Wernly%NULL%1,                    Antonio%Artigas%NULL%1,                    Bernardo Bollen%Pinto%NULL%1,                    Joerg C.%Schefold%NULL%1,                    Dorota%Studzińska%NULL%1,                    Michael%Joannidis%NULL%1,                    Sandra%Oeyen%NULL%1,                    Brian%Marsh%NULL%1,                    Finn H.%Andersen%NULL%1,                    Rui%Moreno%NULL%1,                    Maurizio%Cecconi%NULL%1,                    Christian%Jung%NULL%1,                    Wojciech%Szczeklik%NULL%1,                    NULL%NULL%NULL%0]'
$ws.Range("E5").Value = '[Albert%Prats‐Uribe%NULL%1,                    Marc%Tobed%NULL%1,                    José Miguel%Villacampa%NULL%1,                    Adriana%Agüero%NULL%1,                    Clara%García‐Bastida%NULL%1,                    José Ignacio%Tato%NULL%1,                    Laura%Rodrigáñez%NULL%1,                    Victoria Duque%Holguera%NULL%1,                    Estefanía%Hernández‐García%NULL%1,                    Daniel%Poletti%NULL%1,                    Gabriela%Simonetti%NULL%1,                    Vanessa%Villarraga%NULL%1,                    Carla%Meler‐Claramonte%NULL%1,                    Álvaro%Sánchez Barrueco%NULL%1,                    Carlos%Chiesa‐Estomba%NULL%1,                    Maria%Casasayas%NULL%1,                    Pablo%Parente‐Arias%NULL%1,                    Nieves%Mata‐Castro%NULL%1,                    Jordi%Rello%NULL%1,                    Pedro%Castro%NULL%1,                    Daniel%Prieto‐Alhambra%NULL%1,                    Isabel%Vilaseca%NULL%1,                    Francesc Xavier%Avilés‐Jurado%faviles@clinic.cat%1,                    NULL%NULL%NULL%0]'
$ws.Range("E6").Value = '[T.%Volo%tizianavolo@icloud.com%1,                    P.%Stritoni%NULL%1,                    I.%Battel%NULL%1,                    B.%Zennaro%NULL%1,                    F.%Lazzari%NULL%1,                    M.%Bellin%NULL%1,                    L.%Michieletto%NULL%1,                    G.%Spinato%NULL%1,                    C.%Busatto%NULL%1,                    D.%Politi%NULL%1,                    R.%Spinato%NULL%1]'
$ws.Range("E7").Value = '[Luis F.%Angel%xref no email%1,        Nancy E.%Amoroso%xref no email%1,        Samaan%Rafeq%xref no email%0,        Brian%Mitzman%xref no email%1,        Ronald%Goldenberg%xref no email%1,        Saketh Palasamudram%Shekar%xref no email%1,        Andrea B.%Troxel%xref no email%1,        Yan%Zhang%xref no email%0,        Stephanie H%Chang%xref no email%1,        Paul%Kwak%xref no email%1,        Milan R.%Amin%xref no email%0,        Kimberly%Sureau%xref no email%0,        Heidi B.%Nafday%xref no email%1,        Sarun%Thomas%xref no email%1,        Zachary%Kon%xref no email%1,        Philip M.%Sommer%xref no email%1,        Leopoldo N.%Segal%xref no email%1,        William H.%Moore%xref no email%1,        Robert%Cerfolio%xref no email%1]'
$ws.Range("E8").Value = '[Denise%Battaglini%NULL%1,                    Francesco%Missale%NULL%1,                    Irene%Schiavetti%NULL%1,                    Marta%Filauro%NULL%1,                    Francesca%Iannuzzi%NULL%1,                    Alessandro%Ascoli%NULL%1,                    Alberto%Bertazzoli%NULL%1,                    Federico%Pascucci%NULL%1,                    Salvatore%Grasso%NULL%1,                    Francesco%Murgolo%NULL%1,                    Simone%Binda%NULL%1,                    Davide%Maraggia%NULL%1,                    Giorgia%Montrucchio%NULL%1,                    Gabriele%Sales%NULL%1,                    Giuseppe%Pascarella%NULL%1,                    Felice Eugenio%Agrò%NULL%1,                    Gaia%Faccio%NULL%1,                    Sandra%Ferraris%NULL%1,                    Savino%Spadaro%NULL%1,                    Giulia%Falò%NULL%1,                    Nadia%Mereto%NULL%1,                    Alessandro%Uva%NULL%1,                    Jessica Giuseppina%Maugeri%NULL%1,                    Bellissima%Agrippino%NULL%1,                    Maria%Vargas%NULL%1,                    Giuseppe%Servillo%NULL%1,                    Chiara%Robba%NULL%1,                    Lorenzo%Ball%NULL%1,                    Francesco%Mora%NULL%1,                    Alessio%Signori%NULL%1,                    Antoni%Torres%NULL%1,                    Daniele Roberto%Giacobbe%NULL%1,                    Antonio%Vena%NULL%1,                    Matteo%Bassetti%NULL%1,                    Giorgio%Peretti%NULL%1,                    Patricia R. M.%Rocco%NULL%1,                    Paolo%Pelosi%NULL%1,                    Jihad%Mallat%NULL%1]'
$ws.Range("E9").Value = '[Omar%Breik%NULL%1,                    Paul%Nankivell%NULL%1,                    Neil%Sharma%NULL%1,                    Mansoor N.%Bangash%NULL%1,                    Camilla%Dawson%NULL%1,                    Matthew%Idle%NULL%1,                    Peter%Isherwood%NULL%1,                    Christopher%Jennings%NULL%1,                    Damian%Keene%NULL%1,                    Mav%Manji%NULL%1,                    Tim%Martin%NULL%1,                    Rob%Moss%NULL%1,                    Nick%Murphy%NULL%1,                    Dhruv%Parekh%NULL%1,                    Sat%Parmar%NULL%1,                    Jaimin%Patel%NULL%1,                    Paul%Pracy%NULL%1,                    Prav%Praveen%NULL%1,                    Carla%Richardson%NULL%1,                    Alex%Richter%NULL%1,                    Rajneesh%Sachdeva%NULL%1,                    Adrian%Shields%NULL%1,                    Somiah%Siddiq%NULL%1,                    Simon%Smart%NULL%1,                    Laura%Tasker%NULL%1]'
$ws.Range("E10").Value = '[Toshiki%Kuno%kuno-toshiki@hotmail.co.jp%1,                    Yoshihisa%Miyamoto%NULL%1,                    Matsuo%So%NULL%1,                    Masao%Iwagami%NULL%1,                    Mai%Takahashi%NULL%1,                    Natalia N.%Egorova%NULL%1]'
$ws.Range("E11").Value = '[Kamran%Mahmood%NULL%1,                    George Z.%Cheng%NULL%1,                    Keriann%Van Nostrand%NULL%1,                    Samira%Shojaee%NULL%1,                    Max T.%Wayne%NULL%1,                    Matthew%Abbott%NULL%1,                    Darrell%Nettlow%NULL%1,                    Alice%Parish%NULL%1,                    Cynthia L.%Green%NULL%1,                    Javeryah%Safi%NULL%1,                    Michael J.%Brenner%NULL%1,                    Jose%De Cardenas%NULL%1]'
$ws.Range("E12").Value = '[A%Takhar%NULL%1,                    C%Tornari%NULL%1,                    N%Amin%NULL%1,                    D%Wyncoll%NULL%1,                    S%Tricklebank%NULL%1,                    A%Arora%NULL%1,                    I%Ahmad%NULL%1,                    R%Simo%NULL%1,                    P%Surda%NULL%1]'
$ws.Range("E13").Value = '[Yun%Tang%NULL%1,                    Yongran%Wu%NULL%1,                    Fangfang%Zhu%NULL%1,                    Xiaobo%Yang%NULL%1,                    Chaolin%Huang%NULL%1,                    Guo%Hou%NULL%1,                    Wenhao%Xu%NULL%1,                    Ming%Hu%NULL%1,                    Lu%Zhang%NULL%1,                    Aiguo%Cheng%NULL%1,                    Zhengqin%Xu%NULL%1,                    Boyi%Liu%NULL%1,                    Song%Hu%NULL%1,                    Guochao%Zhu%NULL%1,                    Xuepeng%Fan%NULL%1,                    Xijing%Zhang%NULL%1,                    Yadong%Yang%NULL%1,                    Huibin%Feng%NULL%1,                    Lixia%Yu%NULL%1,                    Bing%Wang%NULL%1,                    Zhiqiang%Li%NULL%1,                    Yong%Peng%NULL%1,                    Zubo%Shen%NULL%1,                    Shouzhi%Fu%NULL%1,                    Yaqi%Ouyang%NULL%1,                    Jiqian%Xu%NULL%1,                    Xiaojing%Zou%NULL%1,                    Minghao%Fang%NULL%1,                    Zhui%Yu%NULL%1,                    Bo%Hu%NULL%1,                    You%Shang%NULL%1]'
